# Generate Report for Handoff
#
# The row for e2e\1d3eb6dd-6952-4a28-b9c4-63947b8ebd7c.md has just been
# handed off for translation. Update its Status from "In Translation" to
# "Ready for handoff" on every sheet that tracks it (Overview, zh-cn,
# de-de), refresh the "Latest Handoff" timestamps for the zh-cn and de-de
# locales, and refresh the handback-staleness hash embedded in the Error
# Detail message for both locales.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is 1d3eb6dd-6952-4a28-b9c4-63947b8ebd7c.md ---
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2017-02-22 08:24:33"

# --- zh-cn sheet: row 3 is 1d3eb6dd-6952-4a28-b9c4-63947b8ebd7c.md ---
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2017-02-22 08:24:16"

$zhcnError = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/7fcf7e2a75201f3fb833729c3cb5244c140dc7bc/e2e/1d3eb6dd-6952-4a28-b9c4-63947b8ebd7c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/90894c2042f210f3f44d3e93013d54039f9baecd/e2e/1d3eb6dd-6952-4a28-b9c4-63947b8ebd7c.md."
$zhcn.Range("R3").Value = $zhcnError

# --- de-de sheet: row 3 is 1d3eb6dd-6952-4a28-b9c4-63947b8ebd7c.md ---
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2017-02-22 08:24:33"

$dedeError = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/7fcf7e2a75201f3fb833729c3cb5244c140dc7bc/e2e/1d3eb6dd-6952-4a28-b9c4-63947b8ebd7c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/90894c2042f210f3f44d3e93013d54039f9baecd/e2e/1d3eb6dd-6952-4a28-b9c4-63947b8ebd7c.md."
$dede.Range("R3").Value = $dedeError
